$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52:127 down to 53:128.
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new data record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this table, so reuse them.
$ws.Range("A52").Value = 3
$ws.Range("B52").Value = "Femacal de La Calera"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = "09/28/2021"
$ws.Range("E52").Value = 5
$ws.Range("F52").Value = 100112010
$ws.Range("G52").Value = "Achicoria"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 140
$ws.Range("K52").Value = 5000
$ws.Range("L52").Value = 5500
$ws.Range("M52").Value = 5286
$ws.Range("N52").Value = "`$/caja 16 unidades"
$ws.Range("O52").Value = "Provincia de Quillota"
$ws.Range("P52").Value = 330
$ws.Range("Q52").Value = 16
$ws.Range("R52").Value = "Hortaliza"
